$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove column A (the old per-row index/counter column). Excel shifts
# columns B:F left by one, turning the old B:F header/data into the new A:E.
$ws.Columns("A").Delete()
